$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value = 518.55554
$ws.Range("I11").Value = 518.55554
$ws.Range("K11").Value = 518.55554
$ws.Range("M11").Value = -378.55554

# row 96
$ws.Range("H96").Value = 1337
$ws.Range("I96").Value = 949.25
$ws.Range("K96").Value = 2847.75
$ws.Range("M96").Value = -1474.75

# row 132
$ws.Range("H132").Value = 1722.8857
$ws.Range("I132").Value = 1468.7142
$ws.Range("K132").Value = 4406.142599999999
$ws.Range("M132").Value = -1876.142599999999

# row 137
$ws.Range("H137").Value = 1268.3422
$ws.Range("I137").Value = 1003.6071
$ws.Range("K137").Value = 3010.8213
$ws.Range("M137").Value = -460.8212999999996

# row 138
$ws.Range("H138").Value = 5684344.5
$ws.Range("I138").Value = 1263
$ws.Range("J138").Value = 8067572
$ws.Range("K138").Value = 3789
$ws.Range("L138").Value = 24202716
$ws.Range("M138").Value = 1351
$ws.Range("N138").Value = -24212996


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3618.4285
$ws.Range("I2").Value = 3065.8
$ws.Range("K2").Value = 3065.8
$ws.Range("M2").Value = -2952.8

# row 88
$ws.Range("H88").Value = 1199.3334
$ws.Range("J88").Value = 969.2
$ws.Range("L88").Value = 969.2
$ws.Range("N88").Value = -1781.2

# row 91
$ws.Range("H91").Value = 1199.3334
$ws.Range("J91").Value = 969.2
$ws.Range("L91").Value = 969.2
$ws.Range("N91").Value = -3777.2

# row 110
$ws.Range("H110").Value = 9499
$ws.Range("I110").Value = 9748.909
$ws.Range("K110").Value = 9748.909
$ws.Range("M110").Value = -7703.909

# row 116
$ws.Range("H116").Value = 3618.4285
$ws.Range("I116").Value = 3065.8
$ws.Range("K116").Value = 3065.8
$ws.Range("M116").Value = -771.8000000000002

# row 134
$ws.Range("H134").Value = 111111
$ws.Range("J134").Value = 111111
$ws.Range("L134").Value = 111111
$ws.Range("N134").Value = -121251


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3618.4285
$ws.Range("I3").Value = 3065.8
$ws.Range("K3").Value = 3065.8
$ws.Range("M3").Value = -2951.8

# row 86
$ws.Range("H86").Value = 2085.963
$ws.Range("I86").Value = 2101.6099
$ws.Range("J86").Value = 2036.6154
$ws.Range("K86").Value = 2101.6099
$ws.Range("L86").Value = 2036.6154
$ws.Range("M86").Value = -978.6098999999999
$ws.Range("N86").Value = -4282.6154

# row 89
$ws.Range("H89").Value = 2085.963
$ws.Range("I89").Value = 2101.6099
$ws.Range("J89").Value = 2036.6154
$ws.Range("K89").Value = 10508.0495
$ws.Range("L89").Value = 10183.077
$ws.Range("M89").Value = -4892.049499999999
$ws.Range("N89").Value = -21415.077

# row 105
$ws.Range("H105").Value = 1070.4166
$ws.Range("I105").Value = 1070.4166
$ws.Range("K105").Value = 1070.4166
$ws.Range("M105").Value = 676.5834

# row 107
$ws.Range("H107").Value = 1022.75
$ws.Range("J107").Value = 1323.75
$ws.Range("L107").Value = 1323.75
$ws.Range("N107").Value = -5163.75


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = $null

# row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = $null

# row 94
$ws.Range("H94").Value = 1746.9231
$ws.Range("J94").Value = 2144.25
$ws.Range("L94").Value = 2144.25
$ws.Range("N94").Value = -3046.25

# row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 27827358
$ws.Range("I4").Value = 28600336
$ws.Range("K4").Value = 85801008
$ws.Range("M4").Value = -85800896

# row 5
$ws.Range("H5").Value = 1226.6666
$ws.Range("I5").Value = 990
$ws.Range("K5").Value = 2970
$ws.Range("M5").Value = -2858

# row 135
$ws.Range("H135").Value = 1226.6666
$ws.Range("I135").Value = 990
$ws.Range("K135").Value = 8910
$ws.Range("M135").Value = -6375


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1334.0667
$ws.Range("I97").Value = 1678.2727
$ws.Range("J97").Value = 387.5
$ws.Range("K97").Value = 1678.2727
$ws.Range("L97").Value = 387.5
$ws.Range("M97").Value = -1182.2727
$ws.Range("N97").Value = -1379.5

# row 122
$ws.Range("H122").Value = 2453.5
$ws.Range("I122").Value = 2320.348
$ws.Range("K122").Value = 6961.044
$ws.Range("M122").Value = -4511.044

# row 132
$ws.Range("H132").Value = 2672.7585
$ws.Range("I132").Value = 2056.3914
$ws.Range("K132").Value = 6169.174199999999
$ws.Range("M132").Value = -3639.174199999999


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 1058.1428
$ws.Range("I46").Value = 793.5
$ws.Range("J46").Value = 1164
$ws.Range("K46").Value = 793.5
$ws.Range("L46").Value = 1164
$ws.Range("M46").Value = -605.5
$ws.Range("N46").Value = -1540

# row 48
$ws.Range("H48").Value = 26899
$ws.Range("I48").Value = 26899
$ws.Range("K48").Value = 26899
$ws.Range("M48").Value = -26238

# row 55
$ws.Range("H55").Value = 102.045456
$ws.Range("I55").Value = 109.210526
$ws.Range("K55").Value = 109.210526
$ws.Range("M55").Value = 63.789474

# row 61
$ws.Range("H61").Value = 3653.16
$ws.Range("J61").Value = 4818.3335
$ws.Range("L61").Value = 4818.3335
$ws.Range("N61").Value = -5222.3335

# row 113
$ws.Range("H113").Value = 3653.16
$ws.Range("J113").Value = 4818.3335
$ws.Range("L113").Value = 4818.3335
$ws.Range("N113").Value = -9158.333500000001

# row 132
$ws.Range("H132").Value = 3544.6
$ws.Range("I132").Value = 3368.9048
$ws.Range("K132").Value = 10106.7144
$ws.Range("M132").Value = -7576.714399999999


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 7776.1113
$ws.Range("I62").Value = 4765
$ws.Range("J62").Value = 8378.333000000001
$ws.Range("K62").Value = 4765
$ws.Range("L62").Value = 8378.333000000001
$ws.Range("M62").Value = -4141
$ws.Range("N62").Value = -9626.333000000001

# row 65
$ws.Range("H65").Value = 7776.1113
$ws.Range("I65").Value = 4765
$ws.Range("J65").Value = 8378.333000000001
$ws.Range("K65").Value = 23825
$ws.Range("L65").Value = 41891.665
$ws.Range("M65").Value = -20705
$ws.Range("N65").Value = -48131.665

# row 100
$ws.Range("H100").Value = 628.13635
$ws.Range("I100").Value = 589.82355
$ws.Range("J100").Value = 758.4
$ws.Range("K100").Value = 1179.6471
$ws.Range("L100").Value = 1516.8
$ws.Range("M100").Value = -638.6470999999999
$ws.Range("N100").Value = -2598.8

